# Completed data extraction and data merge
# - Strip the redundant " County" suffix from the county names in the
#   "Sheet1" extract tab (rows 3-101 first, then row 2 / Alamance last,
#   matching the order the source data was reconciled in).
# - Re-point the view/selection on both tabs to where the edit finished.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Sheet1")
$wsOrig = $wb.Worksheets.Item("orig")

# Strip " County" from the county names, rows 3 through 101 first ...
for ($r = 3; $r -le 101; $r++) {
    $cell = $wsData.Cells.Item($r, 1)
    $name = $cell.Value()
    $cell.Value = ($name -replace " County$", "")
}

# ... then fix up row 2 (Alamance County) last.
$firstCell = $wsData.Cells.Item(2, 1)
$firstName = $firstCell.Value()
$firstCell.Value = ($firstName -replace " County$", "")

# Leave the "orig" summary sheet's selection parked back at the top
# data row inside the frozen pane.
$wsOrig.Activate()
$wsOrig.Range("A5:D5").Select()

# Finish on the "Sheet1" extract tab, selection on the last data row.
$wsData.Activate()
$wsData.Range("A102:D102").Select()
